{"js": "const replacements = [\n  [\"121\u00d72=242\", \"465\u00d74=1860\"],\n  [\"829\u00d76=4974\", \"405\u00d75=2025\"],\n  [\"703\u00d75=3515\", \"151\u00d77=1057\"],\n  [\"985\u00d78=7880\", \"541\u00d73=1623\"],\n  [\"952\u00d72=1904\", \"446\u00d75=2230\"],\n  [\"392\u00d79=3528\", \"102\u00d74=408\"],\n  [\"416\u00d72=832\", \"519\u00d73=1557\"],\n  [\"438\u00d79=3942\", \"415\u00d76=2490\"],\n  [\"708\u00d78=5664\", \"304\u00d77=2128\"],\n  [\"341\u00d73=1023\", \"403\u00d75=2015\"],\n  [\"323\u00d79=2907\", \"679\u00d79=6111\"],\n  [\"511\u00d74=2044\", \"612\u00d76=3672\"],\n  [\"904\u00d77=6328\", \"106\u00d76=636\"],\n  [\"998\u00d75=4990\", \"565\u00d75=2825\"],\n  [\"947\u00d75=4735\", \"120\u00d72=240\"],\n  [\"865\u00d79=7785\", \"359\u00d76=2154\"],\n  [\"742\u00d79=6678\", \"588\u00d74=2352\"],\n  [\"858\u00d75=4290\", \"909\u00d79=8181\"],\n  [\"639\u00d74=2556\", \"794\u00d76=4764\"],\n  [\"877\u00d79=7893\", \"531\u00d77=3717\"],\n  [\"201\u00d72=402\", \"830\u00d76=4980\"],\n  [\"366\u00d73=1098\", \"150\u00d78=1200\"],\n  [\"491\u00d77=3437\", \"838\u00d78=6704\"],\n  [\"883\u00d74=3532\", \"866\u00d75=4330\"],\n  [\"605\u00d77=4235\", \"803\u00d72=1606\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"121\u00d72=242\", \"465\u00d74=1860\")\n    ,@(\"829\u00d76=4974\", \"405\u00d75=2025\")\n    ,@(\"703\u00d75=3515\", \"151\u00d77=1057\")\n    ,@(\"985\u00d78=7880\", \"541\u00d73=1623\")\n    ,@(\"952\u00d72=1904\", \"446\u00d75=2230\")\n    ,@(\"392\u00d79=3528\", \"102\u00d74=408\")\n    ,@(\"416\u00d72=832\", \"519\u00d73=1557\")\n    ,@(\"438\u00d79=3942\", \"415\u00d76=2490\")\n    ,@(\"708\u00d78=5664\", \"304\u00d77=2128\")\n    ,@(\"341\u00d73=1023\", \"403\u00d75=2015\")\n    ,@(\"323\u00d79=2907\", \"679\u00d79=6111\")\n    ,@(\"511\u00d74=2044\", \"612\u00d76=3672\")\n    ,@(\"904\u00d77=6328\", \"106\u00d76=636\")\n    ,@(\"998\u00d75=4990\", \"565\u00d75=2825\")\n    ,@(\"947\u00d75=4735\", \"120\u00d72=240\")\n    ,@(\"865\u00d79=7785\", \"359\u00d76=2154\")\n    ,@(\"742\u00d79=6678\", \"588\u00d74=2352\")\n    ,@(\"858\u00d75=4290\", \"909\u00d79=8181\")\n    ,@(\"639\u00d74=2556\", \"794\u00d76=4764\")\n    ,@(\"877\u00d79=7893\", \"531\u00d77=3717\")\n    ,@(\"201\u00d72=402\", \"830\u00d76=4980\")\n    ,@(\"366\u00d73=1098\", \"150\u00d78=1200\")\n    ,@(\"491\u00d77=3437\", \"838\u00d78=6704\")\n    ,@(\"883\u00d74=3532\", \"866\u00d75=4330\")\n    ,@(\"605\u00d77=4235\", \"803\u00d72=1606\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
